$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44574
$ws.Range("K2").Value = 'Black Amber'
$ws.Range("L2").Value = 'Primera'
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("S2").Value = 1028

# Row 3
$ws.Range("D3").Value = 44706
$ws.Range("L3").Value = 'Segunda'
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("S3").Value = 861

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("S4").Value = 1139

# Row 5
$ws.Range("D5").Value = 44596
$ws.Range("K5").Value = 'Black Amber'
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 861

# Row 6
$ws.Range("D6").Value = 44285
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("S6").Value = 806

# Row 8
$ws.Range("D8").Value = 44580
$ws.Range("M8").Value = 270
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("Q8").Value = '$/bandeja 18 kilos granel'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 1083

# Row 9
$ws.Range("D9").Value = 44238
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 806

# Row 10
$ws.Range("D10").Value = 44238
$ws.Range("K10").Value = 'Fortuna'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 806

# Row 11
$ws.Range("D11").Value = 44278
$ws.Range("K11").Value = 'Angeleno'
$ws.Range("L11").Value = 'Primera'
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 15500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 861

# Row 12
$ws.Range("D12").Value = 44614
$ws.Range("K12").Value = 'Angeleno'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 19000
$ws.Range("P12").Value = 18500
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 1028

# Row 13
$ws.Range("D13").Value = 44245
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 14500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("S13").Value = 806

# Row 14
$ws.Range("D14").Value = 44314
$ws.Range("K14").Value = 'Angeleno'
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("S14").Value = 806

# Row 16
$ws.Range("D16").Value = 44587
$ws.Range("L16").Value = 'Segunda'
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 15500
$ws.Range("S16").Value = 861

# Row 17
$ws.Range("D17").Value = 44944
$ws.Range("K17").Value = 'Larry Ann'
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 7500
$ws.Range("S17").Value = 417

# Row 18
$ws.Range("D18").Value = 44628
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 270
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15500
$ws.Range("S18").Value = 861

# Row 19
$ws.Range("D19").Value = 44239
$ws.Range("L19").Value = 'Primera'
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 15500
$ws.Range("S19").Value = 861

# Row 20
$ws.Range("D20").Value = 44650
$ws.Range("K20").Value = 'Angeleno'
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 17500
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'
$ws.Range("S20").Value = 972

# Row 21
$ws.Range("D21").Value = 44243
$ws.Range("K21").Value = 'Black Amber'
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 14500
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("S21").Value = 806

# Row 22
$ws.Range("D22").Value = 44229
$ws.Range("K22").Value = 'Fortuna'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 14500
$ws.Range("S22").Value = 806

# Row 23
$ws.Range("D23").Value = 44175
$ws.Range("K23").Value = 'Angeleno'
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 21000
$ws.Range("O23").Value = 22000
$ws.Range("P23").Value = 21500
$ws.Range("S23").Value = 1194

# Row 24
$ws.Range("D24").Value = 44169
$ws.Range("K24").Value = 'Angeleno'
$ws.Range("L24").Value = 'Tercera'
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 24500
$ws.Range("S24").Value = 1361
